$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in rows 2 and 3 before removing row 5
$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5

# Remove row 5 (theta_threshold_range), shifting pie_threshold_range up to row 5
$ws.Rows("5:5").Delete()

# Update the (now) row 5 values for pie_threshold_range
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 25

# Update selection to match the saved state
$ws.Range("C5").Select()
